# Update the cover sheet's attendance table: fill in the first blank row
# (the row right after the "24/2/23" entry) with the new session's
# date, time, duration and roles.

$d = $word.ActiveDocument

# The attendance table is the second table in the document (the first is
# the Weighting / Student ID table).
$t = $d.Tables.Item(2)

# Locate the first row below the "24/2/23" row whose first cell is empty -
# that is the row this edit is meant to populate.
$targetRow = 0
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $cellText = $t.Cell($i, 1).Range.Text
    # Cell ranges include trailing cell-mark control characters (CR + BEL)
    # that aren't real content - strip those before checking for blankness.
    $trimmed = ($cellText -replace "[\x07\x0D]", "").Trim()
    if ($trimmed -eq "") {
        $targetRow = $i
        break
    }
}

$t.Cell($targetRow, 1).Range.Text = "28/2/23"
$t.Cell($targetRow, 2).Range.Text = "14:25"
$t.Cell($targetRow, 3).Range.Text = "1h"
$t.Cell($targetRow, 4).Range.Text = "Observer"
$t.Cell($targetRow, 5).Range.Text = "Driver "
